$d = $word.ActiveDocument

# 1) Remove the redundant "και " before "την αποδοτικότητά"
$d.Content.Find.Execute(
    "βελτιώσουν και την αποδοτικότητά", $true, $false, $false, $false, $false,
    $true, 1, $false, "βελτιώσουν την αποδοτικότητά", 2)

# 2) Insert " φύσεως" right after "οικονομικής" (before " αλλά μπορεί")
$d.Content.Find.Execute(
    "είναι οικονομικής αλλά", $true, $false, $false, $false, $false,
    $true, 1, $false, "είναι οικονομικής φύσεως αλλά", 2)
